# Refresh the cryptocurrency price/volume table (coinranking.com snapshot).
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
# All data cells in B:E are stored as literal text in the source workbook, so
# every write below forces text formatting first (and resets the style back to
# "Normal" afterwards) to avoid Excel auto-converting numeric-looking strings
# such as "0.9996" or "27.036.78" into actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# row => @{ Coin = ...; Link = ...; Price = ...; Volume = ... } (only the keys that changed)
$updates = @{
    2 = @{ Price = '27.036.78'; Volume = '  +5.65%  ' }
    3 = @{ Price = '1.878.68' }
    4 = @{ Price = '0.9996'; Volume = '  -0.17%  ' }
    5 = @{ Price = '281.94'; Volume = '  +3.28%  ' }
    6 = @{ Price = '0.9997'; Volume = '  -0.10%  ' }
    7 = @{ Price = '0.5270'; Volume = '  +5.15%  ' }
    8 = @{ Price = '0.3529'; Volume = '  +0.97%  ' }
    9 = @{ Coin = 'OKB'; Link = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; Price = '45.35'; Volume = '  +3.58%  ' }
    10 = @{ Coin = 'Dogecoin'; Link = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; Price = '0.07027'; Volume = '  +6.88%  ' }
    11 = @{ Coin = 'Solana'; Link = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; Price = '20.35'; Volume = '  +2.59%  ' }
    12 = @{ Coin = 'Polygon'; Link = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; Price = '0.8169'; Volume = '  -1.79%  ' }
    13 = @{ Coin = 'TRON'; Link = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Price = '0.07780'; Volume = '  +0.19%  ' }
    14 = @{ Coin = 'WrappedEther'; Link = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Price = '1.876.16'; Volume = '  +4.02%  ' }
    15 = @{ Coin = 'Polkadot'; Link = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Price = '5.210'; Volume = '  +3.43%  ' }
    16 = @{ Coin = 'Litecoin'; Link = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Price = '90.56'; Volume = '  +3.86%  ' }
    17 = @{ Coin = 'BinanceUSD'; Link = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Price = '0.9992'; Volume = '  -0.15%  ' }
    18 = @{ Coin = 'Avalanche'; Link = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Price = '14.58'; Volume = '  +5.27%  ' }
    19 = @{ Coin = 'ShibaInu'; Link = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Price = '0.000008174'; Volume = '  +3.16%  ' }
    20 = @{ Coin = 'Dai'; Link = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; Price = '1.0000'; Volume = '  -0.11%  ' }
    21 = @{ Coin = 'WrappedBTC'; Link = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; Price = '27.063.94'; Volume = '  +5.47%  ' }
    22 = @{ Coin = 'WrappedliquidstakedEther2.0'; Link = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; Price = '2.103.54'; Volume = '  +3.39%  ' }
    23 = @{ Coin = 'Uniswap'; Link = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Price = '4.764'; Volume = '  +1.39%  ' }
    24 = @{ Coin = 'Cosmos'; Link = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Price = '10.18'; Volume = '  +2.15%  ' }
    25 = @{ Coin = 'Chainlink'; Link = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Price = '6.225'; Volume = '  +3.15%  ' }
    26 = @{ Coin = 'LidoDAOToken'; Link = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Price = '2.382'; Volume = '  +13.28%  ' }
    27 = @{ Coin = 'Monero'; Link = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Price = '146.19'; Volume = '  +3.15%  ' }
    28 = @{ Coin = 'EthereumClassic'; Link = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; Price = '17.57'; Volume = '  +4.14%  ' }
    29 = @{ Coin = 'Toncoin'; Link = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Price = '1.680'; Volume = '  +1.77%  ' }
    30 = @{ Coin = 'BitcoinCash'; Link = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; Price = '113.03'; Volume = '  +4.69%  ' }
    31 = @{ Coin = 'InternetComputer(DFINITY)'; Link = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Price = '4.381'; Volume = '  +1.84%  ' }
    32 = @{ Coin = 'Filecoin'; Link = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Price = '4.378'; Volume = '  +4.85%  ' }
    33 = @{ Coin = 'Stellar'; Link = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Price = '0.08890'; Volume = '  +1.51%  ' }
    34 = @{ Coin = 'Hedera'; Link = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; Price = '0.04894'; Volume = '  +2.60%  ' }
    35 = @{ Coin = 'ARBITRUM'; Link = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Price = '1.175'; Volume = '  +4.22%  ' }
    36 = @{ Coin = 'ImmutableX'; Link = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Price = '0.7420'; Volume = '  +3.48%  ' }
    37 = @{ Coin = 'HuobiToken'; Link = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Price = '2.883'; Volume = '  +0.17%  ' }
    38 = @{ Coin = 'MXToken'; Link = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Price = '3.293'; Volume = '  +9.21%  ' }
    39 = @{ Coin = 'RenderToken'; Link = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Price = '2.409'; Volume = '  +7.13%  ' }
    40 = @{ Coin = 'TheSandbox'; Link = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Price = '0.5298'; Volume = '  +3.54%  ' }
    41 = @{ Coin = 'VeChain'; Link = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Price = '0.01881'; Volume = '  +1.51%  ' }
    42 = @{ Coin = 'TrustWalletToken'; Link = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Price = '0.9804'; Volume = '  +4.05%  ' }
    43 = @{ Coin = 'Quant'; Link = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Price = '117.08'; Volume = '  +3.28%  ' }
    44 = @{ Coin = 'FraxShare'; Link = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Price = '6.313'; Volume = '  +3.08%  ' }
    45 = @{ Coin = 'Aptos'; Link = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Price = '8.196'; Volume = '  +3.07%  ' }
    46 = @{ Coin = 'PaxDollar'; Link = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; Price = '0.9990'; Volume = '  -0.13%  ' }
    47 = @{ Coin = 'Decentraland'; Link = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; Price = '0.4596'; Volume = '  +1.51%  ' }
    48 = @{ Coin = 'Algorand'; Link = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Price = '0.1368'; Volume = '  -0.19%  ' }
    49 = @{ Coin = 'EnergySwap'; Link = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Price = '9.448'; Volume = '  +2.72%  ' }
    50 = @{ Coin = 'Elrond'; Link = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'; Price = '36.74'; Volume = '  +2.37%  ' }
    51 = @{ Coin = 'NEARProtocol'; Link = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Price = '1.520'; Volume = '  +2.90%  ' }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    if ($rowData.ContainsKey("Coin"))   { Set-TextValue $ws.Cells.Item($rowNum, 2) $rowData.Coin }
    if ($rowData.ContainsKey("Link"))   { Set-TextValue $ws.Cells.Item($rowNum, 3) $rowData.Link }
    if ($rowData.ContainsKey("Price"))  { Set-TextValue $ws.Cells.Item($rowNum, 4) $rowData.Price }
    if ($rowData.ContainsKey("Volume")) { Set-TextValue $ws.Cells.Item($rowNum, 5) $rowData.Volume }
}
